$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-17: Name, Definition, Found in Document (count), Percentage
$data = @(
    @{ Row = 2;  B = "ORG";         C = "ORG: Companies, agencies, institutions, etc.";                       D = 416; E = 0.5042424242424243 },
    @{ Row = 3;  B = "CARDINAL";    C = "CARDINAL: Numerals that do not fall under another type";             D = 130; E = 0.1575757575757576 },
    @{ Row = 4;  B = "GPE";         C = "GPE: Countries, cities, states";                                     D = 73;  E = 0.08848484848484849 },
    @{ Row = 5;  B = "PRODUCT";     C = "PRODUCT: Objects, vehicles, foods, etc. (not services)";             D = 43;  E = 0.05212121212121212 },
    @{ Row = 6;  B = "LAW";         C = "LAW: Named documents made into laws.";                               D = 31;  E = 0.03757575757575757 },
    @{ Row = 7;  B = "ORDINAL";     C = "ORDINAL: ""first"", ""second"", etc.";                               D = 26;  E = 0.03151515151515152 },
    @{ Row = 8;  B = "DATE";        C = "DATE: Absolute or relative dates or periods";                        D = 26;  E = 0.03151515151515152 },
    @{ Row = 9;  B = "PERSON";      C = "PERSON: People, including fictional";                                D = 25;  E = 0.0303030303030303 },
    @{ Row = 10; B = "LOC";         C = "LOC: Non-GPE locations, mountain ranges, bodies of water";           D = 13;  E = 0.01575757575757576 },
    @{ Row = 11; B = "MONEY";       C = "MONEY: Monetary values, including unit";                             D = 11;  E = 0.01333333333333333 },
    @{ Row = 12; B = "NORP";        C = "NORP: Nationalities or religious or political groups";               D = 9;   E = 0.01090909090909091 },
    @{ Row = 13; B = "PERCENT";     C = "PERCENT: Percentage, including ""%""";                                D = 8;   E = 0.009696969696969697 },
    @{ Row = 14; B = "WORK_OF_ART"; C = "WORK_OF_ART: Titles of books, songs, etc.";                          D = 7;   E = 0.008484848484848486 },
    @{ Row = 15; B = "EVENT";       C = "EVENT: Named hurricanes, battles, wars, sports events, etc.";        D = 5;   E = 0.006060606060606061 },
    @{ Row = 16; B = "FAC";         C = "FAC: Buildings, airports, highways, bridges, etc.";                  D = 1;   E = 0.001212121212121212 },
    @{ Row = 17; B = "TIME";        C = "TIME: Times smaller than a day";                                     D = 1;   E = 0.001212121212121212 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
}
